$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows above row 514, shifting the existing rows 514:537 down to 517:540
$ws.Rows.Item(514).Insert()
$ws.Rows.Item(514).Insert()
$ws.Rows.Item(514).Insert()

# New row 514: Especial
$ws.Cells.Item(514, 1).Value = 7
$ws.Cells.Item(514, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(514, 3).Value = "Ñuble"
$ws.Cells.Item(514, 4).Value = 45041
$ws.Cells.Item(514, 5).Value = 16
$ws.Cells.Item(514, 6).Value = "Fruta"
$ws.Cells.Item(514, 7).Value = 100101
$ws.Cells.Item(514, 8).Value = "Berries"
$ws.Cells.Item(514, 9).Value = 100112025
$ws.Cells.Item(514, 10).Value = "Frutilla"
$ws.Cells.Item(514, 11).Value = "Sin especificar"
$ws.Cells.Item(514, 12).Value = "Especial"
$ws.Cells.Item(514, 13).Value = 50
$ws.Cells.Item(514, 14).Value = 8000
$ws.Cells.Item(514, 15).Value = 8000
$ws.Cells.Item(514, 16).Value = 8000
$ws.Cells.Item(514, 17).Value = "`$/caja 7 kilos"
$ws.Cells.Item(514, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(514, 19).Value = 1143
$ws.Cells.Item(514, 20).Value = 7

# New row 515: Primera
$ws.Cells.Item(515, 1).Value = 7
$ws.Cells.Item(515, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(515, 3).Value = "Ñuble"
$ws.Cells.Item(515, 4).Value = 45041
$ws.Cells.Item(515, 5).Value = 16
$ws.Cells.Item(515, 6).Value = "Fruta"
$ws.Cells.Item(515, 7).Value = 100101
$ws.Cells.Item(515, 8).Value = "Berries"
$ws.Cells.Item(515, 9).Value = 100112025
$ws.Cells.Item(515, 10).Value = "Frutilla"
$ws.Cells.Item(515, 11).Value = "Sin especificar"
$ws.Cells.Item(515, 12).Value = "Primera"
$ws.Cells.Item(515, 13).Value = 50
$ws.Cells.Item(515, 14).Value = 7000
$ws.Cells.Item(515, 15).Value = 7000
$ws.Cells.Item(515, 16).Value = 7000
$ws.Cells.Item(515, 17).Value = "`$/caja 7 kilos"
$ws.Cells.Item(515, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(515, 19).Value = 1000
$ws.Cells.Item(515, 20).Value = 7

# New row 516: Segunda
$ws.Cells.Item(516, 1).Value = 7
$ws.Cells.Item(516, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(516, 3).Value = "Ñuble"
$ws.Cells.Item(516, 4).Value = 45041
$ws.Cells.Item(516, 5).Value = 16
$ws.Cells.Item(516, 6).Value = "Fruta"
$ws.Cells.Item(516, 7).Value = 100101
$ws.Cells.Item(516, 8).Value = "Berries"
$ws.Cells.Item(516, 9).Value = 100112025
$ws.Cells.Item(516, 10).Value = "Frutilla"
$ws.Cells.Item(516, 11).Value = "Sin especificar"
$ws.Cells.Item(516, 12).Value = "Segunda"
$ws.Cells.Item(516, 13).Value = 30
$ws.Cells.Item(516, 14).Value = 6000
$ws.Cells.Item(516, 15).Value = 6000
$ws.Cells.Item(516, 16).Value = 6000
$ws.Cells.Item(516, 17).Value = "`$/caja 7 kilos"
$ws.Cells.Item(516, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(516, 19).Value = 857
$ws.Cells.Item(516, 20).Value = 7
